$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the same date number-format (style index used by the existing
#     START_DATE/END_DATE columns, e.g. E3/F3) to the new date cells so they
#     land on the identical style index instead of minting a new one. ---
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E495:F499").PasteSpecial(-4122) | Out-Null
$ws.Range("E501:F505").PasteSpecial(-4122) | Out-Null
$ws.Range("E507:F511").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New simulated-patient rows (494-511), three patients (1568-1570) each
#     contributing a PER / D_ERA / D_ERA / VIS / D_EXP / D_EXP block. ---
# Row 494
$ws.Range("A494").Value = "PER"
$ws.Range("B494").Value = 1568
$ws.Range("C494").Value = 1960
$ws.Range("K494").Value = 1
$ws.Range("L494").Value = 1

# Row 495
$ws.Range("A495").Value = "D_ERA"
$ws.Range("B495").Value = 1568
$ws.Range("D495").Value = 1398937
$ws.Range("E495").Value = 39491
$ws.Range("F495").Value = 39565

# Row 496
$ws.Range("A496").Value = "D_ERA"
$ws.Range("B496").Value = 1568
$ws.Range("D496").Value = 902427
$ws.Range("E496").Value = 39491
$ws.Range("F496").Value = 39565

# Row 497
$ws.Range("A497").Value = "VIS"
$ws.Range("B497").Value = 1568
$ws.Range("E497").Value = 39491
$ws.Range("F497").Value = 39565
$ws.Range("H497").Value = 9201

# Row 498
$ws.Range("A498").Value = "D_EXP"
$ws.Range("B498").Value = 1568
$ws.Range("D498").Value = 40223504
$ws.Range("E498").Value = 39491
$ws.Range("F498").Value = 39565
$ws.Range("I498").Value = 30
$ws.Range("J498").Value = 30
$ws.Range("M498").Value = "2 times daily"
$ws.Range("N498").Value = "null"
$ws.Range("O498").Value = "null"

# Row 499
$ws.Range("A499").Value = "D_EXP"
$ws.Range("B499").Value = 1568
$ws.Range("D499").Value = 1594707
$ws.Range("E499").Value = 39491
$ws.Range("F499").Value = 39565
$ws.Range("I499").Value = 30
$ws.Range("J499").Value = 30
$ws.Range("M499").Value = "2 times daily"
$ws.Range("N499").Value = "null"
$ws.Range("O499").Value = "null"

# Row 500
$ws.Range("A500").Value = "PER"
$ws.Range("B500").Value = 1569
$ws.Range("C500").Value = 1966
$ws.Range("K500").Value = 1
$ws.Range("L500").Value = 1

# Row 501
$ws.Range("A501").Value = "D_ERA"
$ws.Range("B501").Value = 1569
$ws.Range("D501").Value = 1398937
$ws.Range("E501").Value = 39491
$ws.Range("F501").Value = 39565

# Row 502
$ws.Range("A502").Value = "D_ERA"
$ws.Range("B502").Value = 1569
$ws.Range("D502").Value = 902427
$ws.Range("E502").Value = 39491
$ws.Range("F502").Value = 39565

# Row 503
$ws.Range("A503").Value = "VIS"
$ws.Range("B503").Value = 1569
$ws.Range("E503").Value = 39491
$ws.Range("F503").Value = 39565
$ws.Range("H503").Value = 9201

# Row 504
$ws.Range("A504").Value = "D_EXP"
$ws.Range("B504").Value = 1569
$ws.Range("D504").Value = 19079775
$ws.Range("E504").Value = 39491
$ws.Range("F504").Value = 39565
$ws.Range("I504").Value = 20
$ws.Range("J504").Value = 30
$ws.Range("M504").Value = "4 times daily"
$ws.Range("N504").Value = "null"
$ws.Range("O504").Value = "null"

# Row 505
$ws.Range("A505").Value = "D_EXP"
$ws.Range("B505").Value = 1569
$ws.Range("D505").Value = 902489
$ws.Range("E505").Value = 39491
$ws.Range("F505").Value = 39565
$ws.Range("I505").Value = 5
$ws.Range("J505").Value = 30
$ws.Range("M505").Value = "4 times daily"
$ws.Range("N505").Value = "null"
$ws.Range("O505").Value = "null"

# Row 506
$ws.Range("A506").Value = "PER"
$ws.Range("B506").Value = 1570
$ws.Range("C506").Value = 1950
$ws.Range("K506").Value = 1
$ws.Range("L506").Value = 1

# Row 507
$ws.Range("A507").Value = "D_ERA"
$ws.Range("B507").Value = 1570
$ws.Range("D507").Value = 1398937
$ws.Range("E507").Value = 39491
$ws.Range("F507").Value = 39565

# Row 508
$ws.Range("A508").Value = "D_ERA"
$ws.Range("B508").Value = 1570
$ws.Range("D508").Value = 950370
$ws.Range("E508").Value = 39491
$ws.Range("F508").Value = 39565

# Row 509
$ws.Range("A509").Value = "VIS"
$ws.Range("B509").Value = 1570
$ws.Range("E509").Value = 39491
$ws.Range("F509").Value = 39565
$ws.Range("H509").Value = 9201

# Row 510
$ws.Range("A510").Value = "D_EXP"
$ws.Range("B510").Value = 1570
$ws.Range("D510").Value = 40223506
$ws.Range("E510").Value = 39491
$ws.Range("F510").Value = 39565
$ws.Range("I510").Value = 10
$ws.Range("J510").Value = 10
$ws.Range("M510").Value = "Daily"
$ws.Range("N510").Value = "null"
$ws.Range("O510").Value = "null"

# Row 511
$ws.Range("A511").Value = "D_EXP"
$ws.Range("B511").Value = 1570
$ws.Range("D511").Value = 43219718
$ws.Range("E511").Value = 39491
$ws.Range("F511").Value = 39565
$ws.Range("I511").Value = 30
$ws.Range("J511").Value = 30
$ws.Range("M511").Value = "2 times daily"
$ws.Range("N511").Value = "null"
$ws.Range("O511").Value = "null"

# --- Restore the selection to where the author last left off editing. ---
$ws.Range("L507").Select() | Out-Null
